$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2024-12-19 00:26:37", -0.1194592325112899, -0.001759815370801996, 0.008409047742303097)
    ,@("2024-12-19 00:26:38", -0.1195568323260522, -0.001750704892921996, 0.008372349253418971)
    ,@("2024-12-19 00:26:39", -0.1189354468387319, -0.001725499237453996, 0.008208920913059286)
    ,@("2024-12-19 00:26:40", -0.118844353678287, -0.001878960175965996, 0.008932152307996766)
    ,@("2024-12-19 00:26:41", -0.1191013665238279, -0.001566875694809995, 0.007464681456993709)
    ,@("2024-12-19 00:26:42", -0.1177219558085199, -0.001607417321375996, 0.007569132434914971)
    ,@("2024-12-19 00:26:43", -0.1183628612587927, -0.001585957084591996, 0.007508736734638471)
    ,@("2024-12-19 00:26:44", -0.1178260622775997, -0.001591676440149997, 0.00750163869451607)
    ,@("2024-12-19 00:26:45", -0.1186686740117148, -0.001540607150255995, 0.00731287230775382)
)

$startRow = 191
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}

"Done adding rows 191-199"